# Weekly fruit/vegetable price update ("Fruta / hortaliza, semanal").
#
# Two new weekly price records (for "Ají" / Inferno, date 44508) are
# inserted at the top of the data block (rows 92-93), pushing the existing
# records (previously rows 92-116) down by two rows (to rows 94-118).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 92 - existing rows 92.. shift down to 94..
$ws.Rows.Item(92).Insert()
$ws.Rows.Item(92).Insert()

# New row 92: Ají, Inferno, Primera - fecha 44508
$ws.Range("A92").Value = 8
$ws.Range("B92").Value = "Terminal La Palmera de La Serena"
$ws.Range("C92").Value = "Coquimbo"
$ws.Range("D92").Value = 44508
$ws.Range("E92").Value = 4
$ws.Range("F92").Value = 100112021
$ws.Range("G92").Value = "Ají"
$ws.Range("H92").Value = "Inferno"
$ws.Range("I92").Value = "Primera"
$ws.Range("J92").Value = 600
$ws.Range("K92").Value = 24000
$ws.Range("L92").Value = 25000
$ws.Range("M92").Value = 24500
$ws.Range("N92").Value = "$/caja 12 kilos"
$ws.Range("O92").Value = "Región de Arica y Parinacota"
$ws.Range("P92").Value = 2042
$ws.Range("Q92").Value = 12
$ws.Range("R92").Value = "Hortaliza"

# New row 93: Ají, Inferno, Segunda - fecha 44508
$ws.Range("A93").Value = 8
$ws.Range("B93").Value = "Terminal La Palmera de La Serena"
$ws.Range("C93").Value = "Coquimbo"
$ws.Range("D93").Value = 44508
$ws.Range("E93").Value = 4
$ws.Range("F93").Value = 100112021
$ws.Range("G93").Value = "Ají"
$ws.Range("H93").Value = "Inferno"
$ws.Range("I93").Value = "Segunda"
$ws.Range("J93").Value = 400
$ws.Range("K93").Value = 14000
$ws.Range("L93").Value = 15000
$ws.Range("M93").Value = 14500
$ws.Range("N93").Value = "$/caja 12 kilos"
$ws.Range("O93").Value = "Región de Arica y Parinacota"
$ws.Range("P93").Value = 1208
$ws.Range("Q93").Value = 12
$ws.Range("R93").Value = "Hortaliza"

Write-Host "Inserted 2 new rows (92-93); dimension now $($ws.UsedRange.Address())"
